$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (not auto-converted to numbers)
foreach ($addr in @("D5","D6","D8","D9","D12","D13","D14","D15","D18","D20","D21","D22","D23","D24","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '97.646.22'
$ws.Range('E2').Value = '  -1.59%  '

$ws.Range('D3').Value = '3.410.57'
$ws.Range('E3').Value = '  +2.86%  '

$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('D5').Value = '256.08'
$ws.Range('E5').Value = '  +0.12%  '

$ws.Range('D6').Value = '653.50'
$ws.Range('E6').Value = '  +4.75%  '

$ws.Range('E7').Value = '  -0.27%  '

$ws.Range('D8').Value = '0.427'
$ws.Range('E8').Value = '  +3.55%  '

$ws.Range('D9').Value = '1.04'
$ws.Range('E9').Value = '  +7.37%  '

$ws.Range('E10').Value = '  -0.05%  '

$ws.Range('D11').Value = '3.409.71'
$ws.Range('E11').Value = '  +2.92%  '

$ws.Range('D12').Value = '0.213'
$ws.Range('E12').Value = '  +6.01%  '

$ws.Range('D13').Value = '42.02'
$ws.Range('E13').Value = '  +6.88%  '

$ws.Range('D14').Value = '6.39'
$ws.Range('E14').Value = '  +16.42%  '

$ws.Range('D15').Value = '0.0000258'
$ws.Range('E15').Value = '  +2.40%  '

$ws.Range('D16').Value = '97.448.80'
$ws.Range('E16').Value = '  -1.52%  '

$ws.Range('D17').Value = '4.039.29'
$ws.Range('E17').Value = '  +2.57%  '

$ws.Range('D18').Value = '8.50'
$ws.Range('E18').Value = '  +35.34%  '

$ws.Range('D19').Value = '3.374.95'
$ws.Range('E19').Value = '  +1.82%  '

$ws.Range('D20').Value = '17.38'
$ws.Range('E20').Value = '  +11.52%  '

$ws.Range('D21').Value = '0.508'
$ws.Range('E21').Value = '  +67.89%  '

$ws.Range('D22').Value = '10.81'
$ws.Range('E22').Value = '  +14.41%  '

$ws.Range('D23').Value = '3.44'
$ws.Range('E23').Value = '  -1.12%  '

$ws.Range('D24').Value = '507.58'
$ws.Range('E24').Value = '  +4.03%  '

$ws.Range('E25').Value = '  +0.40%  '

$ws.Range('D26').Value = '6.06'
$ws.Range('E26').Value = '  +7.17%  '

$ws.Range('D27').Value = '95.82'
$ws.Range('E27').Value = '  +7.62%  '

$ws.Range('D28').Value = '12.58'
$ws.Range('E28').Value = '  +4.85%  '

$ws.Range('B29').Value = 'Hedera'
$ws.Range('C29').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D29').Value = '0.150'
$ws.Range('E29').Value = '  +10.76%  '

$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D30').Value = '11.37'
$ws.Range('E30').Value = '  +10.27%  '

$ws.Range('B31').Value = 'Dai'
$ws.Range('C31').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.06%  '

$ws.Range('B32').Value = 'Cronos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D32').Value = '0.195'
$ws.Range('E32').Value = '  +3.59%  '

$ws.Range('B33').Value = 'Binance-PegBSC-USD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.01%  '

$ws.Range('B34').Value = 'PolygonEcosystemToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D34').Value = '0.565'
$ws.Range('E34').Value = '  +20.87%  '

$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').Value = '29.56'
$ws.Range('E35').Value = '  +6.31%  '

$ws.Range('B36').Value = 'PancakeSwap'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D36').Value = '2.17'
$ws.Range('E36').Value = '  +11.51%  '

$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D37').Value = '7.76'
$ws.Range('E37').Value = '  +7.19%  '

$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = '0.155'
$ws.Range('E38').Value = '  +5.03%  '

$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = '509.69'
$ws.Range('E39').Value = '  +3.66%  '

$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = '1.38'
$ws.Range('E40').Value = '  +11.76%  '

$ws.Range('B41').Value = 'WhiteBITCoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D41').Value = '24.69'
$ws.Range('E41').Value = '  -0.67%  '

$ws.Range('B42').Value = 'ARBITRUM'
$ws.Range('C42').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D42').Value = '0.845'
$ws.Range('E42').Value = '  +7.06%  '

$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.0419'
$ws.Range('E43').Value = '  +26.11%  '

$ws.Range('B44').Value = 'MantraDAO'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D44').Value = '3.64'
$ws.Range('E44').Value = '  -0.47%  '

$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '3.30'
$ws.Range('E45').Value = '  +4.83%  '

$ws.Range('B46').Value = 'Filecoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D46').Value = '5.42'
$ws.Range('E46').Value = '  +14.97%  '

$ws.Range('B47').Value = 'USDe'
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').Value = '  +0.00%  '

$ws.Range('D48').Value = '8.13'
$ws.Range('E48').Value = '  +11.01%  '

$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '2.09'
$ws.Range('E49').Value = '  +6.41%  '

$ws.Range('B50').Value = 'ImmutableX'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D50').Value = '1.57'
$ws.Range('E50').Value = '  +15.26%  '

$ws.Range('B51').Value = 'OKB'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D51').Value = '50.49'
$ws.Range('E51').Value = '  +9.64%  '
